$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (July -> August data correction)
$ws.Range("A2").Value = 45139
$ws.Range("B2").Value = 1045
$ws.Range("C2").Value = 30

# Add row 3 (September data)
$ws.Range("A3").Value = 45170
$ws.Range("B3").Value = 1107
$ws.Range("C3").Value = 18

# Update selection to match diff
$ws.Range("D11").Select()
